# Apply "Added Week 15 simulations" updates to Target Depth Data workbook.
# Updates row 3 (the "R" row) on both the OFF and DEF sheets with new
# simulated totals for Short Att (B), Short Comp (C), Deep Att (D), Deep Comp (E).

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 327
$wsOff.Range("C3").Value = 207
$wsOff.Range("D3").Value = 117
$wsOff.Range("E3").Value = 56

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 430
$wsDef.Range("C3").Value = 281
$wsDef.Range("D3").Value = 110
$wsDef.Range("E3").Value = 41
